$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Style = "Normal"
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '67.622.89'
Set-TextValue $ws.Range('E2') '  -0.79%  '
Set-TextValue $ws.Range('D3') '3.251.31'
Set-TextValue $ws.Range('E3') '  -0.67%  '
Set-TextValue $ws.Range('E4') '  -0.01%  '
Set-TextValue $ws.Range('D5') '580.45'
Set-TextValue $ws.Range('E5') '  -1.19%  '
Set-TextValue $ws.Range('D6') '184.84'
Set-TextValue $ws.Range('E6') '  -0.21%  '
Set-TextValue $ws.Range('E7') '  +0.00%  '
Set-TextValue $ws.Range('E8') '  +0.92%  '
Set-TextValue $ws.Range('D9') '3.248.53'
Set-TextValue $ws.Range('E9') '  -0.72%  '
Set-TextValue $ws.Range('D10') '0.132'
Set-TextValue $ws.Range('E10') '  -2.78%  '
Set-TextValue $ws.Range('D11') '6.58'
Set-TextValue $ws.Range('E11') '  -2.41%  '
Set-TextValue $ws.Range('D12') '0.412'
Set-TextValue $ws.Range('E12') '  -1.22%  '
Set-TextValue $ws.Range('D13') '3.809.36'
Set-TextValue $ws.Range('E13') '  -0.84%  '
Set-TextValue $ws.Range('E14') '  -0.11%  '
Set-TextValue $ws.Range('D15') '27.70'
Set-TextValue $ws.Range('E15') '  -3.20%  '
Set-TextValue $ws.Range('D16') '67.636.86'
Set-TextValue $ws.Range('E16') '  -0.76%  '
Set-TextValue $ws.Range('E17') '  -1.47%  '
Set-TextValue $ws.Range('D18') '3.268.52'
Set-TextValue $ws.Range('E18') '  -0.14%  '
Set-TextValue $ws.Range('E19') '  -1.71%  '
Set-TextValue $ws.Range('D20') '13.58'
Set-TextValue $ws.Range('E20') '  -0.56%  '
Set-TextValue $ws.Range('D21') '394.53'
Set-TextValue $ws.Range('E21') '  +3.14%  '
Set-TextValue $ws.Range('D22') '7.60'
Set-TextValue $ws.Range('E22') '  -1.88%  '
Set-TextValue $ws.Range('E23') '  -0.06%  '
Set-TextValue $ws.Range('D24') '71.46'
Set-TextValue $ws.Range('E24') '  -0.01%  '
Set-TextValue $ws.Range('E25') '  +0.53%  '
Set-TextValue $ws.Range('E26') '  -2.54%  '
Set-TextValue $ws.Range('E27') '  -2.86%  '
Set-TextValue $ws.Range('D28') '9.59'
Set-TextValue $ws.Range('E28') '  -1.97%  '
Set-TextValue $ws.Range('D29') '0.999'
Set-TextValue $ws.Range('E29') '  -0.02%  '
Set-TextValue $ws.Range('E30') '  -2.09%  '
Set-TextValue $ws.Range('D31') '5.54'
Set-TextValue $ws.Range('E31') '  -5.10%  '
Set-TextValue $ws.Range('D32') '22.69'
Set-TextValue $ws.Range('E32') '  -0.99%  '
Set-TextValue $ws.Range('D33') '7.00'
Set-TextValue $ws.Range('E33') '  -2.84%  '
Set-TextValue $ws.Range('E34') '  -1.87%  '
Set-TextValue $ws.Range('E35') '  +0.04%  '
Set-TextValue $ws.Range('D36') '161.82'
Set-TextValue $ws.Range('E36') '  -0.91%  '
Set-TextValue $ws.Range('E37') '  -4.24%  '
Set-TextValue $ws.Range('E38') '  +1.25%  '
Set-TextValue $ws.Range('D39') '26.61'
Set-TextValue $ws.Range('E39') '  -0.03%  '
Set-TextValue $ws.Range('D40') '0.808'
Set-TextValue $ws.Range('E40') '  -3.73%  '
Set-TextValue $ws.Range('E41') '  -1.31%  '
Set-TextValue $ws.Range('E42') '  -4.59%  '
Set-TextValue $ws.Range('D43') '2.49'
Set-TextValue $ws.Range('E43') '  -6.42%  '
Set-TextValue $ws.Range('E44') '  -0.26%  '
Set-TextValue $ws.Range('D45') '40.64'
Set-TextValue $ws.Range('E45') '  -1.43%  '
Set-TextValue $ws.Range('D46') '2.616.70'
Set-TextValue $ws.Range('E46') '  -0.79%  '
Set-TextValue $ws.Range('D47') '24.77'
Set-TextValue $ws.Range('E47') '  -3.37%  '
Set-TextValue $ws.Range('D48') '334.31'
Set-TextValue $ws.Range('E48') '  -2.41%  '
Set-TextValue $ws.Range('D50') '6.36'
Set-TextValue $ws.Range('E50') '  +1.62%  '
Set-TextValue $ws.Range('E51') '  -0.29%  '
